# ============================================================
# Daily Report - October: add three new day blocks (Oct 17-19, 2014)
# rows 76-84 (Fri 10/17), rows 85-93 (Sat 10/18), rows 94-102 (Sun 10/19 = NGHI)
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: replicate row structure/styles via copy of existing day blocks ---
# Block 1 (rows 76-84) and Block 2 (rows 85-93): copy the previous weekday block (67-75)
$ws.Range("B67:C75").Copy($ws.Range("B76:C84"))
$ws.Range("B67:C75").Copy($ws.Range("B85:C93"))
# Block 3 (rows 94-102): copy the NGHI (weekend) block (22-30)
$ws.Range("B22:C30").Copy($ws.Range("B94:C102"))

# --- Step 2: fix the date formulas (copy turned them into static values) ---
$ws.Range("B76").Formula = "=B67+1"
$ws.Range("B85").Formula = "=B76+1"
$ws.Range("B94").Formula = "=B85+1"

# --- Step 2b: the "Vấn đề gặp phải" / "Giải quyết vấn đề:" sections are blank
#     for these new days (copy brought over the old day's leftover text) ---
$ws.Range("C80").ClearContents()
$ws.Range("C82").ClearContents()
$ws.Range("C89").ClearContents()
$ws.Range("C91").ClearContents()

# --- Step 3: write the new day's report text, in first-seen order ---
$ws.Range("C78").Value = "- Công việc 1: Hoàn thành
'- Công Việc 2: Hoàn thành
'- Công việc 3: Hoàn thành"
$ws.Range("C93").Value = "- Làm API chi tiết Công việc, Tìm kiếm"
$ws.Range("C77").Value = "- Công việc 1: Hoàn thiện configure widget
'- Công việc 2: Click item trên listview của widget
'- Công việc 3: Update api công việc và công văn"
$ws.Range("C86").Value = "- Công việc 1: Chuyển id của công việc khi click vào item sang màn hình chi tiết công việc"
$ws.Range("C84").Value = "- Chuyển id của công việc khi click vào item sang màn hình chi tiết công việc"
$ws.Range("C87").Value = "- Công việc 1: Hoàn thành"

# --- Step 4: re-apply the original cell formatting (Value= can reset styles) ---
$ws.Range("C68").Copy()
$ws.Range("C77").PasteSpecial(-4122)
$ws.Range("C69").Copy()
$ws.Range("C78").PasteSpecial(-4122)
$ws.Range("C75").Copy()
$ws.Range("C84").PasteSpecial(-4122)
$ws.Range("C68").Copy()
$ws.Range("C86").PasteSpecial(-4122)
$ws.Range("C69").Copy()
$ws.Range("C87").PasteSpecial(-4122)
$ws.Range("C75").Copy()
$ws.Range("C93").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 5: row heights to match the auto-fit wrapping of the longer text ---
$ws.Rows.Item(76).RowHeight = 18
$ws.Rows.Item(77).RowHeight = 42.75
$ws.Rows.Item(78).RowHeight = 42.75
$ws.Rows.Item(85).RowHeight = 18
$ws.Rows.Item(94).RowHeight = 18
$ws.Rows.Item(95).RowHeight = 28.5
$ws.Rows.Item(96).RowHeight = 28.5

# --- Step 6: update sheet view / selection to reflect the newly added bottom rows ---
$ws.Range("D96").Select()
$excel.ActiveWindow.ScrollRow = 85

Write-Output "Added day blocks for 2014-10-17, 2014-10-18, 2014-10-19 (NGHI)."
